$d = $word.ActiveDocument

$d.Content.Find.Execute("2023-05-16 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-05-17 Wednesday", 2) | Out-Null
$d.Content.Find.Execute("78×61=", $true, $false, $false, $false, $false, $true, 1, $false, "79×27=", 2) | Out-Null
$d.Content.Find.Execute("33×54=", $true, $false, $false, $false, $false, $true, 1, $false, "90×48=", 2) | Out-Null
$d.Content.Find.Execute("29×17=", $true, $false, $false, $false, $false, $true, 1, $false, "66×17=", 2) | Out-Null
$d.Content.Find.Execute("53×17=", $true, $false, $false, $false, $false, $true, 1, $false, "96×37=", 2) | Out-Null
$d.Content.Find.Execute("57×48=", $true, $false, $false, $false, $false, $true, 1, $false, "82×48=", 2) | Out-Null
$d.Content.Find.Execute("25×74=", $true, $false, $false, $false, $false, $true, 1, $false, "83×13=", 2) | Out-Null
$d.Content.Find.Execute("98×60=", $true, $false, $false, $false, $false, $true, 1, $false, "12×81=", 2) | Out-Null
$d.Content.Find.Execute("94×19=", $true, $false, $false, $false, $false, $true, 1, $false, "91×91=", 2) | Out-Null
$d.Content.Find.Execute("44×95=", $true, $false, $false, $false, $false, $true, 1, $false, "34×52=", 2) | Out-Null
$d.Content.Find.Execute("45×41=", $true, $false, $false, $false, $false, $true, 1, $false, "63×85=", 2) | Out-Null
$d.Content.Find.Execute("40×12=", $true, $false, $false, $false, $false, $true, 1, $false, "77×68=", 2) | Out-Null
$d.Content.Find.Execute("31×12=", $true, $false, $false, $false, $false, $true, 1, $false, "94×75=", 2) | Out-Null
$d.Content.Find.Execute("14×24=", $true, $false, $false, $false, $false, $true, 1, $false, "18×29=", 2) | Out-Null
$d.Content.Find.Execute("72×46=", $true, $false, $false, $false, $false, $true, 1, $false, "79×85=", 2) | Out-Null
$d.Content.Find.Execute("22×75=", $true, $false, $false, $false, $false, $true, 1, $false, "86×37=", 2) | Out-Null
$d.Content.Find.Execute("17×37=", $true, $false, $false, $false, $false, $true, 1, $false, "40×42=", 2) | Out-Null
$d.Content.Find.Execute("65×59=", $true, $false, $false, $false, $false, $true, 1, $false, "26×29=", 2) | Out-Null
$d.Content.Find.Execute("23×38=", $true, $false, $false, $false, $false, $true, 1, $false, "17×80=", 2) | Out-Null
$d.Content.Find.Execute("47×29=", $true, $false, $false, $false, $false, $true, 1, $false, "77×37=", 2) | Out-Null
$d.Content.Find.Execute("67×45=", $true, $false, $false, $false, $false, $true, 1, $false, "92×61=", 2) | Out-Null
$d.Content.Find.Execute("84×57=", $true, $false, $false, $false, $false, $true, 1, $false, "86×63=", 2) | Out-Null
$d.Content.Find.Execute("90×81=", $true, $false, $false, $false, $false, $true, 1, $false, "16×96=", 2) | Out-Null
$d.Content.Find.Execute("66×79=", $true, $false, $false, $false, $false, $true, 1, $false, "70×19=", 2) | Out-Null
$d.Content.Find.Execute("18×81=", $true, $false, $false, $false, $false, $true, 1, $false, "19×53=", 2) | Out-Null
$d.Content.Find.Execute("84×29=", $true, $false, $false, $false, $false, $true, 1, $false, "17×46=", 2) | Out-Null
$d.Content.Find.Execute("61×85=", $true, $false, $false, $false, $false, $true, 1, $false, "100×65=", 2) | Out-Null
$d.Content.Find.Execute("57×57=", $true, $false, $false, $false, $false, $true, 1, $false, "78×95=", 2) | Out-Null
$d.Content.Find.Execute("33×10=", $true, $false, $false, $false, $false, $true, 1, $false, "30×81=", 2) | Out-Null
$d.Content.Find.Execute("33×81=", $true, $false, $false, $false, $false, $true, 1, $false, "43×77=", 2) | Out-Null
$d.Content.Find.Execute("32×73=", $true, $false, $false, $false, $false, $true, 1, $false, "67×38=", 2) | Out-Null
$d.Content.Find.Execute("96×62=", $true, $false, $false, $false, $false, $true, 1, $false, "33×23=", 2) | Out-Null
$d.Content.Find.Execute("100×76=", $true, $false, $false, $false, $false, $true, 1, $false, "95×23=", 2) | Out-Null
$d.Content.Find.Execute("21×29=", $true, $false, $false, $false, $false, $true, 1, $false, "39×55=", 2) | Out-Null
$d.Content.Find.Execute("93×95=", $true, $false, $false, $false, $false, $true, 1, $false, "44×24=", 2) | Out-Null
$d.Content.Find.Execute("70×76=", $true, $false, $false, $false, $false, $true, 1, $false, "54×91=", 2) | Out-Null
$d.Content.Find.Execute("59×41=", $true, $false, $false, $false, $false, $true, 1, $false, "23×12=", 2) | Out-Null
$d.Content.Find.Execute("31×73=", $true, $false, $false, $false, $false, $true, 1, $false, "60×20=", 2) | Out-Null
$d.Content.Find.Execute("23×50=", $true, $false, $false, $false, $false, $true, 1, $false, "25×91=", 2) | Out-Null
$d.Content.Find.Execute("47×52=", $true, $false, $false, $false, $false, $true, 1, $false, "23×55=", 2) | Out-Null
$d.Content.Find.Execute("38×95=", $true, $false, $false, $false, $false, $true, 1, $false, "32×68=", 2) | Out-Null
$d.Content.Find.Execute("42×29=", $true, $false, $false, $false, $false, $true, 1, $false, "36×10=", 2) | Out-Null
$d.Content.Find.Execute("16×70=", $true, $false, $false, $false, $false, $true, 1, $false, "43×98=", 2) | Out-Null
$d.Content.Find.Execute("83×15=", $true, $false, $false, $false, $false, $true, 1, $false, "96×17=", 2) | Out-Null
$d.Content.Find.Execute("60×24=", $true, $false, $false, $false, $false, $true, 1, $false, "75×43=", 2) | Out-Null
$d.Content.Find.Execute("66×58=", $true, $false, $false, $false, $false, $true, 1, $false, "71×49=", 2) | Out-Null
$d.Content.Find.Execute("40×39=", $true, $false, $false, $false, $false, $true, 1, $false, "30×91=", 2) | Out-Null
$d.Content.Find.Execute("20×52=", $true, $false, $false, $false, $false, $true, 1, $false, "67×19=", 2) | Out-Null
$d.Content.Find.Execute("38×51=", $true, $false, $false, $false, $false, $true, 1, $false, "19×93=", 2) | Out-Null
$d.Content.Find.Execute("20×32=", $true, $false, $false, $false, $false, $true, 1, $false, "94×22=", 2) | Out-Null
$d.Content.Find.Execute("35×20=", $true, $false, $false, $false, $false, $true, 1, $false, "32×38=", 2) | Out-Null
$d.Content.Find.Execute("37×43=", $true, $false, $false, $false, $false, $true, 1, $false, "96×15=", 2) | Out-Null
$d.Content.Find.Execute("42×20=", $true, $false, $false, $false, $false, $true, 1, $false, "37×40=", 2) | Out-Null
$d.Content.Find.Execute("92×27=", $true, $false, $false, $false, $false, $true, 1, $false, "14×39=", 2) | Out-Null
$d.Content.Find.Execute("58×34=", $true, $false, $false, $false, $false, $true, 1, $false, "100×41=", 2) | Out-Null
$d.Content.Find.Execute("93×47=", $true, $false, $false, $false, $false, $true, 1, $false, "76×43=", 2) | Out-Null
$d.Content.Find.Execute("80×74=", $true, $false, $false, $false, $false, $true, 1, $false, "58×98=", 2) | Out-Null
$d.Content.Find.Execute("91×63=", $true, $false, $false, $false, $false, $true, 1, $false, "27×87=", 2) | Out-Null
$d.Content.Find.Execute("70×18=", $true, $false, $false, $false, $false, $true, 1, $false, "23×11=", 2) | Out-Null
$d.Content.Find.Execute("96×78=", $true, $false, $false, $false, $false, $true, 1, $false, "86×34=", 2) | Out-Null
$d.Content.Find.Execute("37×82=", $true, $false, $false, $false, $false, $true, 1, $false, "91×86=", 2) | Out-Null
$d.Content.Find.Execute("82×84=", $true, $false, $false, $false, $false, $true, 1, $false, "32×80=", 2) | Out-Null
$d.Content.Find.Execute("76×97=", $true, $false, $false, $false, $false, $true, 1, $false, "96×14=", 2) | Out-Null
$d.Content.Find.Execute("96×27=", $true, $false, $false, $false, $false, $true, 1, $false, "66×90=", 2) | Out-Null
$d.Content.Find.Execute("61×37=", $true, $false, $false, $false, $false, $true, 1, $false, "72×35=", 2) | Out-Null
$d.Content.Find.Execute("76×30=", $true, $false, $false, $false, $false, $true, 1, $false, "89×58=", 2) | Out-Null
$d.Content.Find.Execute("99×16=", $true, $false, $false, $false, $false, $true, 1, $false, "100×40=", 2) | Out-Null
$d.Content.Find.Execute("87×83=", $true, $false, $false, $false, $false, $true, 1, $false, "35×85=", 2) | Out-Null
$d.Content.Find.Execute("93×54=", $true, $false, $false, $false, $false, $true, 1, $false, "19×31=", 2) | Out-Null
$d.Content.Find.Execute("55×87=", $true, $false, $false, $false, $false, $true, 1, $false, "36×22=", 2) | Out-Null
$d.Content.Find.Execute("12×58=", $true, $false, $false, $false, $false, $true, 1, $false, "63×54=", 2) | Out-Null
$d.Content.Find.Execute("30×95=", $true, $false, $false, $false, $false, $true, 1, $false, "25×62=", 2) | Out-Null
$d.Content.Find.Execute("27×41=", $true, $false, $false, $false, $false, $true, 1, $false, "81×64=", 2) | Out-Null
$d.Content.Find.Execute("54×52=", $true, $false, $false, $false, $false, $true, 1, $false, "58×38=", 2) | Out-Null
$d.Content.Find.Execute("91×69=", $true, $false, $false, $false, $false, $true, 1, $false, "48×58=", 2) | Out-Null
$d.Content.Find.Execute("60×44=", $true, $false, $false, $false, $false, $true, 1, $false, "65×72=", 2) | Out-Null
$d.Content.Find.Execute("88×85=", $true, $false, $false, $false, $false, $true, 1, $false, "58×85=", 2) | Out-Null
$d.Content.Find.Execute("28×52=", $true, $false, $false, $false, $false, $true, 1, $false, "62×57=", 2) | Out-Null
$d.Content.Find.Execute("100×90=", $true, $false, $false, $false, $false, $true, 1, $false, "26×98=", 2) | Out-Null
$d.Content.Find.Execute("95×88=", $true, $false, $false, $false, $false, $true, 1, $false, "93×42=", 2) | Out-Null
$d.Content.Find.Execute("25×57=", $true, $false, $false, $false, $false, $true, 1, $false, "56×50=", 2) | Out-Null
$d.Content.Find.Execute("23×72=", $true, $false, $false, $false, $false, $true, 1, $false, "83×65=", 2) | Out-Null
$d.Content.Find.Execute("36×81=", $true, $false, $false, $false, $false, $true, 1, $false, "12×47=", 2) | Out-Null
$d.Content.Find.Execute("28×12=", $true, $false, $false, $false, $false, $true, 1, $false, "29×81=", 2) | Out-Null
$d.Content.Find.Execute("43×64=", $true, $false, $false, $false, $false, $true, 1, $false, "22×50=", 2) | Out-Null
$d.Content.Find.Execute("43×58=", $true, $false, $false, $false, $false, $true, 1, $false, "14×66=", 2) | Out-Null
$d.Content.Find.Execute("54×10=", $true, $false, $false, $false, $false, $true, 1, $false, "53×24=", 2) | Out-Null
$d.Content.Find.Execute("17×75=", $true, $false, $false, $false, $false, $true, 1, $false, "90×40=", 2) | Out-Null
$d.Content.Find.Execute("60×85=", $true, $false, $false, $false, $false, $true, 1, $false, "73×65=", 2) | Out-Null
$d.Content.Find.Execute("63×71=", $true, $false, $false, $false, $false, $true, 1, $false, "21×74=", 2) | Out-Null
$d.Content.Find.Execute("53×25=", $true, $false, $false, $false, $false, $true, 1, $false, "58×79=", 2) | Out-Null
$d.Content.Find.Execute("92×89=", $true, $false, $false, $false, $false, $true, 1, $false, "91×73=", 2) | Out-Null
$d.Content.Find.Execute("34×24=", $true, $false, $false, $false, $false, $true, 1, $false, "52×71=", 2) | Out-Null
$d.Content.Find.Execute("56×98=", $true, $false, $false, $false, $false, $true, 1, $false, "17×74=", 2) | Out-Null
$d.Content.Find.Execute("23×88=", $true, $false, $false, $false, $false, $true, 1, $false, "51×22=", 2) | Out-Null
$d.Content.Find.Execute("97×53=", $true, $false, $false, $false, $false, $true, 1, $false, "54×53=", 2) | Out-Null
$d.Content.Find.Execute("11×15=", $true, $false, $false, $false, $false, $true, 1, $false, "30×60=", 2) | Out-Null
$d.Content.Find.Execute("43×16=", $true, $false, $false, $false, $false, $true, 1, $false, "44×34=", 2) | Out-Null
$d.Content.Find.Execute("89×54=", $true, $false, $false, $false, $false, $true, 1, $false, "24×32=", 2) | Out-Null
$d.Content.Find.Execute("96×49=", $true, $false, $false, $false, $false, $true, 1, $false, "63×60=", 2) | Out-Null
$d.Content.Find.Execute("29×21=", $true, $false, $false, $false, $false, $true, 1, $false, "12×83=", 2) | Out-Null

Write-Output "Replacements complete"
